# "Changed long term data files"
# The Thursday/Friday columns (D & E) were removed from the top of the
# roster and their rosters appended as additional rows (10-15) under the
# Monday/Tuesday/Wednesday columns (A/B/C). A couple of cells that used to
# carry the styled (Calibri) look in column E are left behind as empty but
# still-styled cells at E6:E10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the old "Thursday 1600" (D) / "Friday 1600" (E) rosters before
# --- we start overwriting cells ---------------------------------------
$d2 = $ws.Range("D2").Value2
$d3 = $ws.Range("D3").Value2
$d4 = $ws.Range("D4").Value2
$d5 = $ws.Range("D5").Value2
$d6 = $ws.Range("D6").Value2
$d7 = $ws.Range("D7").Value2
$d8 = $ws.Range("D8").Value2
$d9 = $ws.Range("D9").Value2

$e2 = $ws.Range("E2").Value2
$e3 = $ws.Range("E3").Value2
$e4 = $ws.Range("E4").Value2
$e5 = $ws.Range("E5").Value2
$e6 = $ws.Range("E6").Value2
$e7 = $ws.Range("E7").Value2
$e8 = $ws.Range("E8").Value2
$e9 = $ws.Range("E9").Value2
$e10 = $ws.Range("E10").Value2

# --- drop the "Thursday 1600" / "Friday 1600" header cells and the whole
# --- column D, then clear out column E content (formatting on E6:E10 is
# --- restored afterwards) ----------------------------------------------
$ws.Range("D1:D10").ClearContents()
$ws.Range("E1:E10").ClearContents()

# --- write the relocated rosters into the new rows 10-15 ---------------
$ws.Range("A10").Value = $d2
$ws.Range("B10").Value = $d6
$ws.Range("C10").Value = $e2

$ws.Range("A11").Value = $d3
$ws.Range("B11").Value = $d7
$ws.Range("C11").Value = $e3

$ws.Range("A12").Value = $d4
$ws.Range("B12").Value = $d8
$ws.Range("C12").Value = $e4

$ws.Range("A13").Value = $d5
$ws.Range("B13").Value = $d9
$ws.Range("C13").Value = $e5

$ws.Range("A14").Value = $e6
$ws.Range("B14").Value = $e8
$ws.Range("C14").Value = $e10

$ws.Range("A15").Value = $e7
$ws.Range("B15").Value = $e9

# --- re-apply the Calibri-look formatting that used to live on E6:E9 to
# --- the new styled cells (A14:C14, A15:B15) and restore the leftover
# --- empty-but-styled cells at E6:E10 -----------------------------------
$ws.Range("E6").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A15:B15").PasteSpecial(-4122)
$ws.Range("E6:E10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- clear stray content so E6:E10 are empty again (format only) -------
$ws.Range("E6:E10").ClearContents()

# --- final cosmetics: dimension follows automatically from used cells,
# --- move the active selection to where the author left off ------------
$ws.Range("D10").Select()
